# CORRECCIONES DE HORARIO Y VALIDACIÓN DE ASISTENCIA
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("F2").Value = "EA1"

# --- Row 3 ---
$ws.Range("F3").Value = "EA2"

# --- Row 4 ---
$ws.Range("F4").Value = "EB1"

# --- Row 5 ---
$ws.Range("A5").Value = "05GRG"
$ws.Range("D5").Value = "23:00"
$ws.Range("E5").Value = "00:00"
$ws.Range("F5").Value = "EB2"

# --- Row 6 ---
$ws.Range("D6").Value = "00:00"
$ws.Range("E6").Value = "01:00"
$ws.Range("F6").Value = "EB3"

# --- Row 7 ---
$ws.Range("F7").Value = "EB3"

# --- Row 8 ---
$ws.Range("F8").Value = "EC4"

# --- Row 9 ---
$ws.Range("F9").Value = "EC1"

# --- Row 10 ---
$ws.Range("F10").Value = "EC2"

# --- Row 11 ---
$ws.Range("B11").Value = "04IB"

# --- Row 12 : clear out the whole row (only D12/E12 keep their text style, now blank) ---
$ws.Range("A12").ClearContents()
$ws.Range("B12").ClearContents()
$ws.Range("C12").ClearContents()
$ws.Range("D12").ClearContents()
$ws.Range("E12").ClearContents()
$ws.Range("F12").ClearContents()
$ws.Range("G12").ClearContents()

# --- Update the selected cell shown when the workbook was last saved ---
$ws.Range("E6").Select()
